$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '38.158.23'
$ws.Range('E2').Value = '  +3.24%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.059.70'
$ws.Range('E3').Value = '  +3.27%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '230.16'
$ws.Range('E5').Value = '  +2.41%  '
$ws.Range('E6').Value = '  +1.69%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '58.03'
$ws.Range('E7').Value = '  +6.39%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.388'
$ws.Range('E9').Value = '  +3.22%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0808'
$ws.Range('E10').Value = '  +3.47%  '
$ws.Range('E11').Value = '  -0.48%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.363.78'
$ws.Range('E12').Value = '  +3.32%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.63'
$ws.Range('E13').Value = '  +4.22%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.64'
$ws.Range('E14').Value = '  +2.92%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.754'
$ws.Range('E15').Value = '  +2.48%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.29'
$ws.Range('E16').Value = '  +4.10%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.064.31'
$ws.Range('E17').Value = '  +3.19%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '38.057.37'
$ws.Range('E18').Value = '  +3.34%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.17'
$ws.Range('E19').Value = '  +1.58%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '69.92'
$ws.Range('E20').Value = '  +2.04%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0832'
$ws.Range('E21').Value = '  +2.40%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '224.99'
$ws.Range('E22').Value = '  +1.20%  '
$ws.Range('E23').Value = '  +0.05%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.44'
$ws.Range('E24').Value = '  +1.46%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.25'
$ws.Range('E25').Value = '  +3.83%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '167.09'
$ws.Range('E26').Value = '  +0.94%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.32'
$ws.Range('E27').Value = '  +2.55%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.134'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.08'
$ws.Range('E29').Value = '  +2.59%  '
$ws.Range('E30').Value = '  +2.36%  '
$ws.Range('E31').Value = '  +1.90%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.56'
$ws.Range('E32').Value = '  +1.77%  '
$ws.Range('E33').Value = '  +5.42%  '
$ws.Range('E34').Value = '  +1.72%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.98'
$ws.Range('E35').Value = '  +7.15%  '
$ws.Range('E36').Value = '  +2.91%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.09'
$ws.Range('E37').Value = '  +15.62%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.33'
$ws.Range('E38').Value = '  +6.23%  '
$ws.Range('E39').Value = '  +0.04%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0220'
$ws.Range('E40').Value = '  +2.55%  '
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '98.49'
$ws.Range('E41').Value = '  +4.32%  '
$ws.Range('B42').Value = 'InjectiveProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '17.07'
$ws.Range('E42').Value = '  +4.85%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.483.17'
$ws.Range('E43').Value = '  +1.17%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0944'
$ws.Range('E44').Value = '  +3.39%  '
$ws.Range('E45').Value = '  +4.01%  '
$ws.Range('E46').Value = '  +0.82%  '
$ws.Range('E47').Value = '  +18.22%  '
$ws.Range('E48').Value = '  +1.84%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.96'
$ws.Range('E49').Value = '  +2.33%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.12'
$ws.Range('E50').Value = '  -0.20%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.252.64'
$ws.Range('E51').Value = '  +3.36%  '
